$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1089
$ws.Range("F3").Value = 4695
$ws.Range("F4").Value = 611
$ws.Range("F5").Value = 187
$ws.Range("F6").Value = 1870
$ws.Range("F7").Value = 47
$ws.Range("F8").Value = 766
$ws.Range("F11").Value = 418
$ws.Range("F14").Value = 826
$ws.Range("F15").Value = 1848
$ws.Range("F16").Value = 568
$ws.Range("F17").Value = 526
$ws.Range("F19").Value = 199
$ws.Range("F20").Value = 12
$ws.Range("F21").Value = 12
$ws.Range("F22").Value = 1560
$ws.Range("F23").Value = 1198
$ws.Range("F24").Value = 608
$ws.Range("F25").Value = 2522
$ws.Range("F26").Value = 6
$ws.Range("F28").Value = 1597
$ws.Range("F30").Value = 496
$ws.Range("F33").Value = 4300
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F12").Value = 37
$ws.Range("F18").Value = 284
$ws.Range("F25").Value = 52
$ws.Range("F38").Value = 39
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 1350
$ws.Range("F5").Value = 1738
$ws.Range("F6").Value = 1088
$ws.Range("F7").Value = 311
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 1350
$ws.Range("F4").Value = 1738
$ws.Range("F5").Value = 1088
$ws.Range("F7").Value = 1089
$ws.Range("F9").Value = 4695
$ws.Range("F10").Value = 611
$ws.Range("F11").Value = 187
$ws.Range("F12").Value = 1870
$ws.Range("F13").Value = 47
$ws.Range("F14").Value = 766
$ws.Range("F19").Value = 418
$ws.Range("F22").Value = 37
$ws.Range("F25").Value = 826
$ws.Range("F26").Value = 1848
$ws.Range("F27").Value = 568
$ws.Range("F28").Value = 526
$ws.Range("F31").Value = 12
$ws.Range("F34").Value = 284
$ws.Range("F36").Value = 1560
$ws.Range("F37").Value = 1198
$ws.Range("F39").Value = 2522
$ws.Range("F41").Value = 6
$ws.Range("F42").Value = 52
$ws.Range("F45").Value = 1597
$ws.Range("F46").Value = 496
$ws.Range("F49").Value = 4300
